# Update statistics values on the "stats" worksheet to reflect results from
# the merged-framework path-checking change described in the commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stats")

# Map of cell -> new value for the two duplicate result blocks
# (rows 2-6 correspond to "run 0" results, rows 8-12 are the "Average" block).
$updates = @{
    "D2" = 0.0004260065034031868
    "E2" = 0.1802364899776876
    "G2" = 0.007283815648406744
    "H2" = 0.01174637861549854
    "I2" = 0.09016905352473259
    "J2" = 0.06408471148461103
    "K2" = 0.001822048798203468

    "D3" = 0.009474779944866896
    "E3" = 0.144201691262424
    "G3" = 0.005109652411192656
    "H3" = 0.01787005364894867
    "I3" = 0.05712607409805059
    "J3" = 0.05948562314733863
    "K3" = 0.001261336728930473

    "D4" = 0.0146218161098659
    "E4" = 0.2114646527916193
    "G4" = 0.006648097187280655
    "H4" = 0.02546107722446322
    "I4" = 0.08648623572662473
    "J4" = 0.08666269900277257
    "K4" = 0.001597088761627674

    "D5" = 0.0003816564567387104
    "E5" = 0.1162972971796989
    "G5" = 0.004708411172032356
    "H5" = 0.008516036905348301
    "I5" = 0.0561062884517014
    "J5" = 0.04243542021140456
    "K5" = 0.001228276174515486

    "D6" = 0.0144741921685636
    "E6" = 0.880797129124403
    "G6" = 0.0124666839838028
    "H6" = 0.03878881921991706
    "I6" = 0.7028697719797492
    "J6" = 0.1117980950511992
    "K6" = 0.004529315978288651

    "D8" = 0.0004260065034031868
    "E8" = 0.1802364899776876
    "G8" = 0.007283815648406744
    "H8" = 0.01174637861549854
    "I8" = 0.09016905352473259
    "J8" = 0.06408471148461103
    "K8" = 0.001822048798203468

    "D9" = 0.009474779944866896
    "E9" = 0.144201691262424
    "G9" = 0.005109652411192656
    "H9" = 0.01787005364894867
    "I9" = 0.05712607409805059
    "J9" = 0.05948562314733863
    "K9" = 0.001261336728930473

    "D10" = 0.0146218161098659
    "E10" = 0.2114646527916193
    "G10" = 0.006648097187280655
    "H10" = 0.02546107722446322
    "I10" = 0.08648623572662473
    "J10" = 0.08666269900277257
    "K10" = 0.001597088761627674

    "D11" = 0.0003816564567387104
    "E11" = 0.1162972971796989
    "G11" = 0.004708411172032356
    "H11" = 0.008516036905348301
    "I11" = 0.0561062884517014
    "J11" = 0.04243542021140456
    "K11" = 0.001228276174515486

    "D12" = 0.0144741921685636
    "E12" = 0.880797129124403
    "G12" = 0.0124666839838028
    "H12" = 0.03878881921991706
    "I12" = 0.7028697719797492
    "J12" = 0.1117980950511992
    "K12" = 0.004529315978288651
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
